$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.448.02"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "'1.901.79"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'325.18"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "'0.4831"
$ws.Range("E7").Value = "  +3.79%  "
$ws.Range("D8").Value = "'0.4059"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").Value = "'0.08127"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").Value = "'1.006"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "'23.47"
$ws.Range("E11").Value = "  +5.22%  "
$ws.Range("D12").Value = "'1.901.85"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "'5.987"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "'7.060"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "'90.25"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "'0.06744"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "'0.00001037"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "'17.59"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "'29.463.11"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "'5.555"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'11.79"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("D24").Value = "'2.155"
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").Value = "'2.149.82"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("D26").Value = "'153.54"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'20.01"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").Value = "'6.223"
$ws.Range("E28").Value = "  +8.58%  "
$ws.Range("D29").Value = "'2.092"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("D30").Value = "'118.78"
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("D31").Value = "'1.025"
$ws.Range("E31").Value = "  -4.86%  "
$ws.Range("D32").Value = "'0.09542"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("D33").Value = "'5.502"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").Value = "'3.555"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "'1.387"
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("D36").Value = "'0.02260"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "'0.06094"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "'1.170"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").Value = "'7.901"
$ws.Range("E40").Value = "  -6.06%  "
$ws.Range("D41").Value = "'0.1851"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").Value = "'10.28"
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("D43").Value = "'1.285"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").Value = "'2.396"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").Value = "'0.07707"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'12.37"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("D47").Value = "'0.5539"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "'114.94"
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("D50").Value = "'72.50"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("E51").Value = "  +2.14%  "
